$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every date/time in column A (rows 2-97) forward by exactly 1 day,
# keeping the same time-of-day fraction.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Mark a new "on" block in column B for rows 36-39 (payload = 1)
for ($r = 36; $r -le 39; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Update the view: top-left visible cell and active selection
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E38").Select()
